# Update losses and own use subsector labellings
# (renumbering / renaming of the 10_01_xx "own use" subsectors and the
#  10_02 "losses" label on the 9th_EBT_schema sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F21: losses label - fix typo "transmision" -> "transmission"
$ws.Range("F21").Value = "10_02_transmission_and_distribution_losses"

# G41:G54: own-use subsector labels renumbered/renamed
$ws.Range("G41").Value = "10_01_13_pump_storage_plants"
$ws.Range("G42").Value = "10_01_05_coke_ovens"
$ws.Range("G43").Value = "10_01_06_coal_mines"
$ws.Range("G44").Value = "10_01_07_blast_furnaces"
$ws.Range("G45").Value = "10_01_11_oil_refineries"
$ws.Range("G46").Value = "10_01_12_oil_and_gas_extraction"
$ws.Range("G47").Value = "10_01_15_charcoal_production_plants"
$ws.Range("G48").Value = "10_01_18_ccs"
$ws.Range("G49").Value = "10_01_16_gasification_plants_for_biogases"
$ws.Range("G50").Value = "10_01_04_gastoliquids_plants"
$ws.Range("G51").Value = "10_01_09_bkb_pb_plants"
$ws.Range("G52").Value = "10_01_17_nonspecified_own_uses"
$ws.Range("G54").Value = "10_01_10_liquefaction_plants_coal_to_oil"

# Restore the cursor/viewport position recorded in the saved file
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("F22").Select()
